$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Cells whose new text would otherwise be auto-parsed as a number by Excel
# get a temporary "@" (Text) number format while the value is entered, then
# the format is cleared again so the cell keeps the original (General/style 0)
# formatting but retains the literal text value (e.g. "1.0000", not 1).
$forceTextRange = $ws.Range("D5,D6,D7,D8,D9,D10,D11,D13,D14,D15,D16,D18,D19,D20,D22,D23,D24,D25,D26,D29,D30,D33,D35,D36,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50")
$forceTextRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.341.78"
$ws.Range("D3").Value = "1.872.39"
$ws.Range("D5").Value = "243.78"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D7").Value = "0.4690"
$ws.Range("D8").Value = "0.2876"
$ws.Range("D9").Value = "0.06451"
$ws.Range("D10").Value = "22.07"
$ws.Range("D11").Value = "0.07768"
$ws.Range("D12").Value = "1.876.47"
$ws.Range("D13").Value = "95.53"
$ws.Range("D14").Value = "0.7215"
$ws.Range("D15").Value = "5.133"
$ws.Range("D16").Value = "279.44"
$ws.Range("D17").Value = "30.332.20"
$ws.Range("D18").Value = "12.99"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D20").Value = "0.000007453"
$ws.Range("D21").Value = "2.117.76"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D23").Value = "5.244"
$ws.Range("D24").Value = "6.242"
$ws.Range("D25").Value = "163.53"
$ws.Range("D26").Value = "9.062"
$ws.Range("D29").Value = "1.317"
$ws.Range("D30").Value = "0.09575"
$ws.Range("D33").Value = "4.097"
$ws.Range("D35").Value = "1.118"
$ws.Range("D36").Value = "0.6883"
$ws.Range("D39").Value = "2.811"
$ws.Range("D40").Value = "6.219"
$ws.Range("D41").Value = "74.29"
$ws.Range("D42").Value = "1.945"
$ws.Range("D43").Value = "0.4229"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D45").Value = "0.8236"
$ws.Range("D46").Value = "100.81"
$ws.Range("D47").Value = "9.564"
$ws.Range("D48").Value = "35.17"
$ws.Range("D49").Value = "6.925"
$ws.Range("D50").Value = "898.06"

# Clear the temporary text formatting so styling matches the original cells
$forceTextRange.ClearFormats()

# --- Column E (Volume 1h) updates ---
# These values already contain non-numeric characters (%, padding spaces)
# so Excel keeps them as plain text without any special handling needed.
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("E51").Value = "  +0.89%  "
